$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.934942333333333
$ws.Range("H2").Value = 23.804827
$ws.Range("I2").Value = 0.1481951088167925
$ws.Range("J2").Value = 0.1539507021887745
$ws.Range("M2").Value = 23.80409633333333
$ws.Range("N2").Value = 71.41228899999999
$ws.Range("O2").Value = 0.2035379592047788
$ws.Range("P2").Value = 0.2090753787995941
$ws.Range("Q2").Value = 188.8841317021114
$ws.Range("R2").Value = 1699.957185319003
$ws.Range("S2").Value = 0.03016333001270006
$ws.Range("T2").Value = 0.03218730137658153
$ws.Range("G3").Value = 7.934942333333333
$ws.Range("H3").Value = 23.804827
$ws.Range("I3").Value = 0.1481951088167925
$ws.Range("J3").Value = 0.1539507021887745
$ws.Range("O3").Value = 0.3053048738509329
$ws.Range("P3").Value = 0.3136109470643028
$ws.Range("Q3").Value = 283.3242812645937
$ws.Range("R3").Value = 2549.918531381343
$ws.Range("S3").Value = 0.0452446890026361
$ws.Range("T3").Value = 0.04828062551463601
$ws.Range("G4").Value = 7.934942333333333
$ws.Range("H4").Value = 23.804827
$ws.Range("I4").Value = 0.1481951088167925
$ws.Range("J4").Value = 0.1539507021887745
$ws.Range("M4").Value = 33.79564933333334
$ws.Range("N4").Value = 101.386948
$ws.Range("O4").Value = 0.2889711669362822
$ws.Range("P4").Value = 0.2968328680576918
$ws.Range("Q4").Value = 268.1665285775551
$ws.Range("R4").Value = 2413.498757197996
$ws.Range("S4").Value = 0.04282411352903785
$ws.Range("T4").Value = 0.0456976284701895
$ws.Range("G5").Value = 7.934942333333333
$ws.Range("H5").Value = 23.804827
$ws.Range("I5").Value = 0.1481951088167925
$ws.Range("J5").Value = 0.1539507021887745
$ws.Range("M5").Value = 9.2924895
$ws.Range("N5").Value = 18.584979
$ws.Range("O5").Value = 0.07945583492339121
$ws.Range("P5").Value = 0.05441166469831967
$ws.Range("Q5").Value = 73.7353683156055
$ws.Range("R5").Value = 442.412209893633
$ws.Range("S5").Value = 0.01177496610260106
$ws.Range("T5").Value = 0.008376713987566467
$ws.Range("G6").Value = 7.934942333333333
$ws.Range("H6").Value = 23.804827
$ws.Range("I6").Value = 0.1481951088167925
$ws.Range("J6").Value = 0.1539507021887745
$ws.Range("M6").Value = 14.353493
$ws.Range("N6").Value = 43.060479
$ws.Range("O6").Value = 0.122730165084615
$ws.Range("P6").Value = 0.1260691413800917
$ws.Range("Q6").Value = 113.8941392369037
$ws.Range("R6").Value = 1025.047253132133
$ws.Range("S6").Value = 0.01818801016981742
$ws.Range("T6").Value = 0.01940843283980101
$ws.Range("I7").Value = 0.1782715511859743
$ws.Range("J7").Value = 0.1851952517494486
$ws.Range("M7").Value = 23.80409633333333
$ws.Range("N7").Value = 71.41228899999999
$ws.Range("O7").Value = 0.2035379592047788
$ws.Range("P7").Value = 0.2090753787995941
$ws.Range("Q7").Value = 227.218478543576
$ws.Range("R7").Value = 2044.966306892184
$ws.Range("S7").Value = 0.03628502771266347
$ws.Range("T7").Value = 0.03871976741140215
$ws.Range("I8").Value = 0.1782715511859743
$ws.Range("J8").Value = 0.1851952517494486
$ws.Range("O8").Value = 0.3053048738509329
$ws.Range("P8").Value = 0.3136109470643028
$ws.Range("S8").Value = 0.05442717344604402
$ws.Range("T8").Value = 0.05807925829295655
$ws.Range("I9").Value = 0.1782715511859743
$ws.Range("J9").Value = 0.1851952517494486
$ws.Range("M9").Value = 33.79564933333334
$ws.Range("N9").Value = 101.386948
$ws.Range("O9").Value = 0.2889711669362822
$ws.Range("P9").Value = 0.2968328680576918
$ws.Range("Q9").Value = 322.5913689552321
$ws.Range("R9").Value = 2903.322320597088
$ws.Range("S9").Value = 0.05151533817775217
$ws.Range("T9").Value = 0.05497203772745508
$ws.Range("I10").Value = 0.1782715511859743
$ws.Range("J10").Value = 0.1851952517494486
$ws.Range("M10").Value = 9.2924895
$ws.Range("N10").Value = 18.584979
$ws.Range("O10").Value = 0.07945583492339121
$ws.Range("P10").Value = 0.05441166469831967
$ws.Range("Q10").Value = 88.70008323380401
$ws.Range("R10").Value = 532.2004994028241
$ws.Range("S10").Value = 0.01416471494256966
$ws.Range("T10").Value = 0.0100767819419119
$ws.Range("I11").Value = 0.1782715511859743
$ws.Range("J11").Value = 0.1851952517494486
$ws.Range("M11").Value = 14.353493
$ws.Range("N11").Value = 43.060479
$ws.Range("O11").Value = 0.122730165084615
$ws.Range("P11").Value = 0.1260691413800917
$ws.Range("Q11").Value = 137.009143114536
$ws.Range("R11").Value = 1233.082288030824
$ws.Range("S11").Value = 0.02187929690694502
$ws.Range("T11").Value = 0.02334740637572291
$ws.Range("G12").Value = 14.055385
$ws.Range("H12").Value = 42.166155
$ws.Range("I12").Value = 0.2625021357479615
$ws.Range("J12").Value = 0.2726971790574536
$ws.Range("M12").Value = 23.80409633333333
$ws.Range("N12").Value = 71.41228899999999
$ws.Range("O12").Value = 0.2035379592047788
$ws.Range("P12").Value = 0.2090753787995941
$ws.Range("Q12").Value = 334.5757385420883
$ws.Range("R12").Value = 3011.181646878795
$ws.Range("S12").Value = 0.05342914899703589
$ws.Range("T12").Value = 0.05701426600901785
$ws.Range("G13").Value = 14.055385
$ws.Range("H13").Value = 42.166155
$ws.Range("I13").Value = 0.2625021357479615
$ws.Range("J13").Value = 0.2726971790574536
$ws.Range("O13").Value = 0.3053048738509329
$ws.Range("P13").Value = 0.3136109470643028
$ws.Range("Q13").Value = 501.860213437655
$ws.Range("R13").Value = 4516.741920938895
$ws.Range("S13").Value = 0.08014318144013186
$ws.Range("T13").Value = 0.0855208205859718
$ws.Range("G14").Value = 14.055385
$ws.Range("H14").Value = 42.166155
$ws.Range("I14").Value = 0.2625021357479615
$ws.Range("J14").Value = 0.2726971790574536
$ws.Range("M14").Value = 33.79564933333334
$ws.Range("N14").Value = 101.386948
$ws.Range("O14").Value = 0.2889711669362822
$ws.Range("P14").Value = 0.2968328680576918
$ws.Range("Q14").Value = 475.0108627049934
$ws.Range("R14").Value = 4275.097764344941
$ws.Range("S14").Value = 0.07585554849035479
$ws.Range("T14").Value = 0.08094548577086587
$ws.Range("G15").Value = 14.055385
$ws.Range("H15").Value = 42.166155
$ws.Range("I15").Value = 0.2625021357479615
$ws.Range("J15").Value = 0.2726971790574536
$ws.Range("M15").Value = 9.2924895
$ws.Range("N15").Value = 18.584979
$ws.Range("O15").Value = 0.07945583492339121
$ws.Range("P15").Value = 0.05441166469831967
$ws.Range("Q15").Value = 130.6095175309575
$ws.Range("R15").Value = 783.6571051857451
$ws.Range("S15").Value = 0.02085732636502766
$ws.Range("T15").Value = 0.01483790747105181
$ws.Range("G16").Value = 14.055385
$ws.Range("H16").Value = 42.166155
$ws.Range("I16").Value = 0.2625021357479615
$ws.Range("J16").Value = 0.2726971790574536
$ws.Range("M16").Value = 14.353493
$ws.Range("N16").Value = 43.060479
$ws.Range("O16").Value = 0.122730165084615
$ws.Range("P16").Value = 0.1260691413800917
$ws.Range("Q16").Value = 201.743870209805
$ws.Range("R16").Value = 1815.694831888245
$ws.Range("S16").Value = 0.03221693045541132
$ws.Range("T16").Value = 0.03437869922054631
$ws.Range("G17").Value = 6.0053675
$ws.Range("H17").Value = 12.010735
$ws.Range("I17").Value = 0.1121578522894532
$ws.Range("J17").Value = 0.07767588846805276
$ws.Range("M17").Value = 23.80409633333333
$ws.Range("N17").Value = 71.41228899999999
$ws.Range("O17").Value = 0.2035379592047788
$ws.Range("P17").Value = 0.2090753787995941
$ws.Range("Q17").Value = 142.9523464870691
$ws.Range("R17").Value = 857.7140789224148
$ws.Range("S17").Value = 0.02282838036378633
$ws.Range("T17").Value = 0.01624011580505315
$ws.Range("G18").Value = 6.0053675
$ws.Range("H18").Value = 12.010735
$ws.Range("I18").Value = 0.1121578522894532
$ws.Range("J18").Value = 0.07767588846805276
$ws.Range("O18").Value = 0.3053048738509329
$ws.Range("P18").Value = 0.3136109470643028
$ws.Range("Q18").Value = 214.4270694343525
$ws.Range("R18").Value = 1286.562416606115
$ws.Range("S18").Value = 0.03424233894462308
$ws.Range("T18").Value = 0.02436000894652719
$ws.Range("G19").Value = 6.0053675
$ws.Range("H19").Value = 12.010735
$ws.Range("I19").Value = 0.1121578522894532
$ws.Range("J19").Value = 0.07767588846805276
$ws.Range("M19").Value = 33.79564933333334
$ws.Range("N19").Value = 101.386948
$ws.Range("O19").Value = 0.2889711669362822
$ws.Range("P19").Value = 0.2968328680576918
$ws.Range("Q19").Value = 202.9552941477967
$ws.Range("R19").Value = 1217.73176488678
$ws.Range("S19").Value = 0.03241038545715046
$ws.Range("T19").Value = 0.02305675675290149
$ws.Range("G20").Value = 6.0053675
$ws.Range("H20").Value = 12.010735
$ws.Range("I20").Value = 0.1121578522894532
$ws.Range("J20").Value = 0.07767588846805276
$ws.Range("M20").Value = 9.2924895
$ws.Range("N20").Value = 18.584979
$ws.Range("O20").Value = 0.07945583492339121
$ws.Range("P20").Value = 0.05441166469831967
$ws.Range("Q20").Value = 55.80481443739125
$ws.Range("R20").Value = 223.219257749565
$ws.Range("S20").Value = 0.008911595796872888
$ws.Range("T20").Value = 0.004226474398467762
$ws.Range("G21").Value = 6.0053675
$ws.Range("H21").Value = 12.010735
$ws.Range("I21").Value = 0.1121578522894532
$ws.Range("J21").Value = 0.07767588846805276
$ws.Range("M21").Value = 14.353493
$ws.Range("N21").Value = 43.060479
$ws.Range("O21").Value = 0.122730165084615
$ws.Range("P21").Value = 0.1260691413800917
$ws.Range("Q21").Value = 86.1980003736775
$ws.Range("R21").Value = 517.1880022420651
$ws.Range("S21").Value = 0.01376515172702045
$ws.Range("T21").Value = 0.009792532565103181
$ws.Range("G22").Value = 16.00284133333333
$ws.Range("H22").Value = 48.00852399999999
$ws.Range("I22").Value = 0.2988733519598186
$ws.Range("J22").Value = 0.3104809785362705
$ws.Range("M22").Value = 23.80409633333333
$ws.Range("N22").Value = 71.41228899999999
$ws.Range("O22").Value = 0.2035379592047788
$ws.Range("P22").Value = 0.2090753787995941
$ws.Range("Q22").Value = 380.933176705715
$ws.Range("R22").Value = 3428.398590351435
$ws.Range("S22").Value = 0.06083207211859305
$ws.Range("T22").Value = 0.06491392819753941
$ws.Range("G23").Value = 16.00284133333333
$ws.Range("H23").Value = 48.00852399999999
$ws.Range("I23").Value = 0.2988733519598186
$ws.Range("J23").Value = 0.3104809785362705
$ws.Range("O23").Value = 0.3053048738509329
$ws.Range("P23").Value = 0.3136109470643028
$ws.Range("Q23").Value = 571.3959003723907
$ws.Range("R23").Value = 5142.563103351515
$ws.Range("S23").Value = 0.09124749101749789
$ws.Range("T23").Value = 0.09737023372421127
$ws.Range("G24").Value = 16.00284133333333
$ws.Range("H24").Value = 48.00852399999999
$ws.Range("I24").Value = 0.2988733519598186
$ws.Range("J24").Value = 0.3104809785362705
$ws.Range("M24").Value = 33.79564933333334
$ws.Range("N24").Value = 101.386948
$ws.Range("O24").Value = 0.2889711669362822
$ws.Range("P24").Value = 0.2968328680576918
$ws.Range("Q24").Value = 540.8264140383058
$ws.Range("R24").Value = 4867.437726344751
$ws.Range("S24").Value = 0.08636578128198698
$ws.Range("T24").Value = 0.09216095933627982
$ws.Range("G25").Value = 16.00284133333333
$ws.Range("H25").Value = 48.00852399999999
$ws.Range("I25").Value = 0.2988733519598186
$ws.Range("J25").Value = 0.3104809785362705
$ws.Range("M25").Value = 9.2924895
$ws.Range("N25").Value = 18.584979
$ws.Range("O25").Value = 0.07945583492339121
$ws.Range("P25").Value = 0.05441166469831967
$ws.Range("Q25").Value = 148.706235060166
$ws.Range("R25").Value = 892.2374103609959
$ws.Range("S25").Value = 0.02374723171631995
$ws.Range("T25").Value = 0.01689378689932174
$ws.Range("G26").Value = 16.00284133333333
$ws.Range("H26").Value = 48.00852399999999
$ws.Range("I26").Value = 0.2988733519598186
$ws.Range("J26").Value = 0.3104809785362705
$ws.Range("M26").Value = 14.353493
$ws.Range("N26").Value = 43.060479
$ws.Range("O26").Value = 0.122730165084615
$ws.Range("P26").Value = 0.1260691413800917
$ws.Range("Q26").Value = 229.6966710581107
$ws.Range("R26").Value = 2067.270039522996
$ws.Range("S26").Value = 0.03668077582542077
$ws.Range("T26").Value = 0.03914207037891831
